# Bug fixes and GUI layout improvements
#  - Append the new "New Dataset" placeholder row to the Datasets sheet
#  - Fix the Methods sheet's Duval (C++) entry to point at the release build
#  - Make the Methods sheet the active/selected sheet again, with the
#    selection covering the whole used range on each sheet (mirrors what
#    Excel leaves behind after the Datasets_Edit / Methods_Edit figures are
#    closed and the workbook is saved)

$wb = $excel.ActiveWorkbook

$wsMethods  = $wb.Worksheets.Item("Methods")
$wsDatasets = $wb.Worksheets.Item("Datasets")

# --- Datasets sheet: append a new row (13, "New Dataset", placeholder path, " Nothing") ---
$wsDatasets.Cells.Item(14, 1).Value = 13
$wsDatasets.Cells.Item(14, 2).Value = "New Dataset"
$wsDatasets.Cells.Item(14, 3).Value = "..\..\..\..\..\documents\Book1.xlsx"
$wsDatasets.Cells.Item(14, 4).Value = " Nothing"

# --- Methods sheet: "Duval (C++)" now points at the release build, not debug ---
$wsMethods.Cells.Item(10, 3).Value = ".\methods\cpp\release\DGA_Test.exe"

# --- refresh selections to span the used data range on each sheet ---
$wsMethods.Range("A1:D10").Select()
$wsDatasets.Range("A1:D14").Select()

# --- Methods becomes the active/selected tab again ---
$wsMethods.Activate()
$wsMethods.Range("A1:D10").Select()
